$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "category" column (C) values for several transactions
$ws.Range("C2").Value = "food"
$ws.Range("C4").Value = "shopping"
$ws.Range("C5").Value = "salary"
$ws.Range("C6").Value = "other"
$ws.Range("C7").Value = "other"
$ws.Range("C8").Value = "shopping"
$ws.Range("C9").Value = "other"

# Update the view: scroll so column B is the leftmost visible column,
# and move the active selection to D5
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D5").Select()
